$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 used to hold a throw-away numeric placeholder (0) styled bold with a
# border, while the real payload (now reformatted as pretty-printed JSON
# instead of a single-line Python literal) lived in A2 as a shared string.
# Strip A1's old formatting first (while it still holds the placeholder),
# then move the payload into A1, and finally drop the now-empty A2 cell so
# the sheet only contains the single A1 cell again.
$ws.Range("A1").ClearFormats()

$newText = @'
questions = [
    {
        "title": "You want to build a list for a new marketing campaign.  How can you check if a specific contact meets your list criteria before you add a filter?",
        "ques_type": 2,
        "options": [
            "Click the Actions drop-down and then Edit columns (A).",
            "Click the Actions drop-down and then Delete list (B).",
            "Click Test and then select the contact to test (C).",
            "Click Edit filter on the top left menu and then Edit filter (D)."
        ],
        "score": "Click Test and then select the contact to test (C)."
    },
    {
        "title": "Your manager has noticed a backlog of unanswered requests from the chat of the company website, and he gives you the task of providing an answer to all these pending requests.  How can you easily find all requests not yet closed from the Live Chat?",
        "ques_type": 2,
        "options": [
            "Conversations &gt Inbox &gt Filter: Date = Last Year",
            "Conversations &gt Inbox &gt Filter: Status = Open &gt Channel = Live Chat",
            "Conversations &gt Inbox &gt Left Menu &gt Unassigned &gt Channel = Live Chat",
            "Conversations &gt Inbox &gt Left Menu &gt Assigned to me"
        ],
        "score": "Conversations &gt Inbox &gt Filter: Status = Open &gt Channel = Live Chat"
    },
    {
        "title": "You are going through some deals in HubSpot and want to check them by specific close dates. Upon sorting deals based on their \"Close Date,\" which sorting option does HubSpot provide?",
        "ques_type": 15,
        "options": [
            "Today",
            "Most recent",
            "Next week",
            "Deal owner",
            "Last two years"
        ],
        "score": [
            "Today",
            "Next week"
        ]
    },
    {
        "title": "You are managing a HubSpot campaign in which tracking code monitors traffic from all website pages, landing pages, blogs, and any other pages.  How long does it take to start a new session in the campaign when a user is inactive on a page?",
        "ques_type": 2,
        "options": [
            "30 minutes",
            "1 hour",
            "12 hours",
            "1 day"
        ],
        "score": "30 minutes"
    }
]
'@
$ws.Range("A1").Value = $newText

$ws.Range("A2").ClearContents()

# Setting the multi-line payload auto-expands row 1's height; auto-fit it
# back down so the row keeps using the sheet's default height like before.
$ws.Rows.Item(1).AutoFit()
